# Automatische test-sync: 2025-08-03 14:44:50
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Logs" sheet: append row 16 with the new test-mail entry
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(16, 1).Value = "Hebben we EcoPro-700 nog op voorraad?"
$logs.Cells.Item(16, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(16, 3).Value = "Testmail #6: Hebben we EcoPro-700 nog op voorraad?"
$logs.Cells.Item(16, 4).Value = "Inkoop / Bestellingen"
$logs.Cells.Item(16, 5).Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@bedrijf.nl."
$logs.Cells.Item(16, 6).Value = "2025-08-03 14:44:39"
$logs.Cells.Item(16, 7).Value = "Ja"
$logs.Cells.Item(16, 8).Value = "Ja"
$logs.Cells.Item(16, 9).Value = "Nee"
$logs.Cells.Item(16, 10).Value = "Nee"

# Extend the conditional-formatting ranges so they cover the new row too.
$colRanges = @{
    "D" = "D2:D16"
    "G" = "G2:G16"
    "H" = "H2:H16"
    "I" = "I2:I16"
    "J" = "J2:J16"
}

foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range($col + "2:" + $col + "15")
    $newRange = $logs.Range($colRanges[$col])
    $cfs = $oldRange.FormatConditions
    for ($i = 1; $i -le $cfs.Count; $i++) {
        $cfs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 2) "Dashboard" sheet: append row 5 with the new category tally
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(5, 1).Value = "Inkoop / Bestellingen"
$dash.Cells.Item(5, 2).Value = 1

# ---------------------------------------------------------------------
# 3) Chart on the Dashboard sheet: grow the series references to
#    include the newly added row 5.
# ---------------------------------------------------------------------
$chart = $dash.ChartObjects(1).Chart
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES('Dashboard'!B1,'Dashboard'!`$A`$2:`$A`$5,'Dashboard'!`$B`$2:`$B`$5,1)"
